$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set NumberFormat to Text ("@") per-cell before assigning, so numeric-looking
# strings (e.g. "1.00", "1.60") are preserved as text, matching the original
# inlineStr cell type, without touching formatting of untouched cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.413.35'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.718.57'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -6.11%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '502.89'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.75'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.729.18'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.84%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.29%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.36%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.26%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.203.83'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -5.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.538.21'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.63'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.734.90'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.35%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.76'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.94'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -5.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '342.88'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.35%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.75'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -5.18%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.49'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0827'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.12'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.60'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '151.32'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.82%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.42'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.50%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.19'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.31%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.04%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.82'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.27%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -6.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.55'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.190.64'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.09%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.996'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.601'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.99'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -8.15%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.75'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.45%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.44%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.05'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.31%  '
